# Auto-generated from the unified diff: update Price (D) / Volume(1h) (E)
# text cells for the crypto ranking sheet. Values are stored as plain text
# (not numbers/percentages), so we force text via NumberFormat "@" before
# assigning, exactly as real-world automation does to stop Excel's
# automatic number/percentage type inference from kicking in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}

Set-TextValue "D2" "331.67"
Set-TextValue "E2" "0.50%"
Set-TextValue "E3" "2.21%"
Set-TextValue "D4" "5.683"
Set-TextValue "E4" "-3.06%"
Set-TextValue "D5" "0.08077"
Set-TextValue "E5" "-0.64%"
Set-TextValue "D6" "2.041"
Set-TextValue "E6" "3.13%"
Set-TextValue "D7" "8.735"
Set-TextValue "E7" "-0.50%"
Set-TextValue "D8" "4.541"
Set-TextValue "E8" "-1.60%"
Set-TextValue "D10" "0.9235"
Set-TextValue "E10" "-2.38%"
Set-TextValue "E11" "-4.58%"
Set-TextValue "D12" "0.1951"
Set-TextValue "E12" "-2.48%"
Set-TextValue "D13" "8.736"
Set-TextValue "E13" "-3.30%"
Set-TextValue "D14" "0.09510"
Set-TextValue "E14" "1.69%"
Set-TextValue "D15" "0.03754"
Set-TextValue "E15" "8.76%"
Set-TextValue "D16" "0.1054"
Set-TextValue "E16" "9.76%"
Set-TextValue "E17" "-1.89%"
Set-TextValue "D18" "0.006312"
Set-TextValue "E18" "0.14%"
Set-TextValue "E19" "-0.03%"
Set-TextValue "E20" "-1.74%"
Set-TextValue "D21" "0.1419"
Set-TextValue "E21" "0.80%"
Set-TextValue "D22" "0.2656"
Set-TextValue "E22" "10.11%"
Set-TextValue "D23" "0.04434"
Set-TextValue "E23" "-0.15%"
Set-TextValue "D24" "0.001262"
Set-TextValue "E24" "-0.18%"
Set-TextValue "D25" "0.004294"
Set-TextValue "E25" "-2.90%"
Set-TextValue "E26" "13.84%"
Set-TextValue "D39" "0.02871"
Set-TextValue "E39" "16.17%"
Set-TextValue "D40" "0.05493"
Set-TextValue "E40" "4.00%"
Set-TextValue "D41" "0.007771"
Set-TextValue "E41" "3.35%"
Set-TextValue "D42" "0.009944"
Set-TextValue "E42" "9.60%"
Set-TextValue "D43" "0.1420"
Set-TextValue "E43" "-1.02%"
Set-TextValue "D44" "0.002124"
Set-TextValue "E44" "3.49%"
Set-TextValue "D45" "0.01179"
Set-TextValue "E45" "12.14%"
Set-TextValue "D46" "0.00006782"
Set-TextValue "E46" "-0.83%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.07%"
Set-TextValue "D48" "0.002285"
Set-TextValue "E48" "26.74%"
Set-TextValue "D49" "0.003016"
Set-TextValue "E49" "-13.81%"
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.07%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.07%"

Write-Output "Updated 68 cells"
